$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "OTRO CAMBIO MAS"
$ws.Range("D6").Select()
